$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format to preserve exact numeric-looking strings
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.381.77'
$ws.Range("E2").Value = '  +1.94%  '
$ws.Range("D3").Value = '1.826.10'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '313.03'
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("D6").Value = '1.0000'
$ws.Range("D7").Value = '0.4465'
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").Value = '0.3764'
$ws.Range("E8").Value = '  +2.42%  '
$ws.Range("D9").Value = '0.07408'
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("D10").Value = '0.8773'
$ws.Range("E10").Value = '  +2.49%  '
$ws.Range("D11").Value = '20.84'
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("D12").Value = '1.829.82'
$ws.Range("E12").Value = '  +1.32%  '
$ws.Range("D13").Value = '6.708'
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("D14").Value = '5.422'
$ws.Range("E14").Value = '  +2.03%  '
$ws.Range("D15").Value = '92.87'
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("D16").Value = '0.07074'
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").Value = '0.000008822'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D20").Value = '15.07'
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("D21").Value = '27.390.06'
$ws.Range("E21").Value = '  +1.86%  '
$ws.Range("D22").Value = '5.345'
$ws.Range("E22").Value = '  +3.62%  '
$ws.Range("D23").Value = '10.93'
$ws.Range("E23").Value = '  +0.57%  '
$ws.Range("E24").Value = '  -1.65%  '
$ws.Range("D25").Value = '151.10'
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("D26").Value = '2.268'
$ws.Range("E26").Value = '  +3.81%  '
$ws.Range("D27").Value = '18.65'
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("D28").Value = '5.349'
$ws.Range("E28").Value = '  +2.44%  '
$ws.Range("D29").Value = '117.20'
$ws.Range("E29").Value = '  +0.49%  '
$ws.Range("D30").Value = '0.08915'
$ws.Range("E30").Value = '  +1.02%  '
$ws.Range("D31").Value = '0.7949'
$ws.Range("E31").Value = '  +5.51%  '
$ws.Range("D32").Value = '1.196'
$ws.Range("E32").Value = '  +1.77%  '
$ws.Range("D33").Value = '4.551'
$ws.Range("E33").Value = '  +1.95%  '
$ws.Range("E34").Value = '  +0.97%  '
$ws.Range("D35").Value = '0.9997'
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  +1.60%  '
$ws.Range("D37").Value = '0.01979'
$ws.Range("E37").Value = '  +0.53%  '
$ws.Range("D38").Value = '0.05272'
$ws.Range("E38").Value = '  +1.51%  '
$ws.Range("D39").Value = '7.366'
$ws.Range("E39").Value = '  +5.09%  '
$ws.Range("D40").Value = '0.5342'
$ws.Range("E40").Value = '  +0.29%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '2.355'
$ws.Range("E41").Value = '  +18.47%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").Value = '2.875'
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("D43").Value = '0.1702'
$ws.Range("E43").Value = '  +0.67%  '
$ws.Range("D44").Value = '8.679'
$ws.Range("E44").Value = '  +2.59%  '
$ws.Range("E45").Value = '  -1.65%  '
$ws.Range("D46").Value = '10.60'
$ws.Range("E46").Value = '  +0.69%  '
$ws.Range("D47").Value = '105.26'
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("D48").Value = '1.686'
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("E50").Value = '  +0.80%  '
$ws.Range("D51").Value = '66.09'
$ws.Range("E51").Value = '  +5.42%  '

# Reset style on column D back to default (Normal) to avoid leaving a text-forced number format
$dRange.Style = "Normal"

Write-Host "Applied crypto price updates"